$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1930.2778
$ws.Range("I40").Value = 1924.2
$ws.Range("J40").Value = 1944.091
$ws.Range("K40").Value = 1924.2
$ws.Range("L40").Value = 1944.091
$ws.Range("M40").Value = -1749.2
$ws.Range("N40").Value = -2294.091

$ws.Range("H86").Value = 1659.0667
$ws.Range("I86").Value = 1597.2
$ws.Range("J86").Value = 1690
$ws.Range("K86").Value = 1597.2
$ws.Range("L86").Value = 1690
$ws.Range("M86").Value = -474.2
$ws.Range("N86").Value = -3936

$ws.Range("H89").Value = 1659.0667
$ws.Range("I89").Value = 1597.2
$ws.Range("J89").Value = 1690
$ws.Range("K89").Value = 7986
$ws.Range("L89").Value = 8450
$ws.Range("M89").Value = -2370
$ws.Range("N89").Value = -19682

$ws.Range("H112").Value = 16735081
$ws.Range("J112").Value = 25101082
$ws.Range("L112").Value = 75303246
$ws.Range("N112").Value = -75305462

$ws.Range("H141").Value = 533894.1
$ws.Range("I141").Value = 1433
$ws.Range("J141").Value = 779645.4
$ws.Range("K141").Value = 4299
$ws.Range("L141").Value = 2338936.2
$ws.Range("M141").Value = 881
$ws.Range("N141").Value = -2349296.2


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3816.1753
$ws.Range("I32").Value = 3751.7605
$ws.Range("K32").Value = 3751.7605
$ws.Range("M32").Value = -3464.7605

$ws.Range("H61").Value = 3960.4
$ws.Range("I61").Value = 3739.3333
$ws.Range("J61").Value = 4107.778
$ws.Range("K61").Value = 3739.3333
$ws.Range("L61").Value = 4107.778
$ws.Range("M61").Value = -3527.3333
$ws.Range("N61").Value = -4531.778

$ws.Range("H63").Value = 2625
$ws.Range("I63").Value = 2250
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2250
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1564
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 2625
$ws.Range("I66").Value = 2250
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 11250
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -7818
$ws.Range("N66").Value = -21864

$ws.Range("H80").Value = 25615.818
$ws.Range("J80").Value = 25615.818
$ws.Range("L80").Value = 25615.818
$ws.Range("N80").Value = -27611.818

$ws.Range("H83").Value = 25615.818
$ws.Range("J83").Value = 25615.818
$ws.Range("L83").Value = 76847.454
$ws.Range("N83").Value = -86831.454

$ws.Range("H97").Value = 468.17392
$ws.Range("I97").Value = 429.9
$ws.Range("J97").Value = 723.3333
$ws.Range("K97").Value = 429.9
$ws.Range("L97").Value = 723.3333
$ws.Range("M97").Value = 66.10000000000002
$ws.Range("N97").Value = -1715.3333

$ws.Range("H110").Value = 1486.9697
$ws.Range("I110").Value = 625.4211
$ws.Range("J110").Value = 2656.2144
$ws.Range("K110").Value = 625.4211
$ws.Range("L110").Value = 2656.2144
$ws.Range("M110").Value = 1419.5789
$ws.Range("N110").Value = -6746.2144

$ws.Range("H122").Value = 3333.2222
$ws.Range("I122").Value = 2357
$ws.Range("K122").Value = 7071
$ws.Range("M122").Value = -4621

$ws.Range("H134").Value = 31248.5
$ws.Range("J134").Value = 31248.5
$ws.Range("L134").Value = 31248.5
$ws.Range("N134").Value = -41388.5

$ws.Range("H136").Value = 3960.4
$ws.Range("I136").Value = 3739.3333
$ws.Range("J136").Value = 4107.778
$ws.Range("K136").Value = 11217.9999
$ws.Range("L136").Value = 12323.334
$ws.Range("M136").Value = -8667.999899999999
$ws.Range("N136").Value = -17423.334


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4136.533
$ws.Range("I99").Value = 3670.6667
$ws.Range("K99").Value = 3670.6667
$ws.Range("M99").Value = -2172.6667

$ws.Range("H107").Value = 2112
$ws.Range("I107").Value = 1673.5333
$ws.Range("J107").Value = 3756.25
$ws.Range("K107").Value = 1673.5333
$ws.Range("L107").Value = 3756.25
$ws.Range("M107").Value = 246.4666999999999
$ws.Range("N107").Value = -7596.25

$ws.Range("H134").Value = 2980.0454
$ws.Range("I134").Value = 2645.762
$ws.Range("K134").Value = 7937.286
$ws.Range("M134").Value = -5402.286


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1404.8667
$ws.Range("I16").Value = 744.8182
$ws.Range("J16").Value = 3220
$ws.Range("K16").Value = 744.8182
$ws.Range("L16").Value = 3220
$ws.Range("M16").Value = -457.8182
$ws.Range("N16").Value = -3794

$ws.Range("H58").Value = 11908235
$ws.Range("I58").Value = 2281.682
$ws.Range("J58").Value = 25004782
$ws.Range("K58").Value = 2281.682
$ws.Range("L58").Value = 25004782
$ws.Range("M58").Value = -2078.682
$ws.Range("N58").Value = -25005188

$ws.Range("H62").Value = 3305.0833
$ws.Range("I62").Value = 2428.3572
$ws.Range("J62").Value = 4532.5
$ws.Range("K62").Value = 2428.3572
$ws.Range("L62").Value = 4532.5
$ws.Range("M62").Value = -1804.3572
$ws.Range("N62").Value = -5780.5

$ws.Range("H65").Value = 3305.0833
$ws.Range("I65").Value = 2428.3572
$ws.Range("J65").Value = 4532.5
$ws.Range("K65").Value = 12141.786
$ws.Range("L65").Value = 22662.5
$ws.Range("M65").Value = -9021.786
$ws.Range("N65").Value = -28902.5

$ws.Range("H107").Value = 1047.0952
$ws.Range("I107").Value = 643.55554
$ws.Range("J107").Value = 1349.75
$ws.Range("K107").Value = 643.55554
$ws.Range("L107").Value = 1349.75
$ws.Range("M107").Value = 1276.44446
$ws.Range("N107").Value = -5189.75

$ws.Range("H113").Value = 1404.8667
$ws.Range("I113").Value = 744.8182
$ws.Range("J113").Value = 3220
$ws.Range("K113").Value = 744.8182
$ws.Range("L113").Value = 3220
$ws.Range("M113").Value = 1425.1818
$ws.Range("N113").Value = -7560

$ws.Range("H136").Value = 11908235
$ws.Range("I136").Value = 2281.682
$ws.Range("J136").Value = 25004782
$ws.Range("K136").Value = 6845.045999999999
$ws.Range("L136").Value = 75014346
$ws.Range("M136").Value = -4295.045999999999
$ws.Range("N136").Value = -75019446

$ws.Range("H141").Value = 24026.316
$ws.Range("J141").Value = 24026.316
$ws.Range("L141").Value = 24026.316
$ws.Range("N141").Value = -34386.316


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.92857
$ws.Range("I2").Value = 18.222221
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 109.333326
$ws.Range("L2").Value = 390
$ws.Range("M2").Value = 3.666674
$ws.Range("N2").Value = -616

$ws.Range("H3").Value = 7566
$ws.Range("I3").Value = 2832.5
$ws.Range("J3").Value = 26500
$ws.Range("K3").Value = 8497.5
$ws.Range("L3").Value = 79500
$ws.Range("M3").Value = -8385.5
$ws.Range("N3").Value = -79724

$ws.Range("H56").Value = 5362.3076
$ws.Range("I56").Value = 5362.3076
$ws.Range("K56").Value = 5362.3076
$ws.Range("M56").Value = -4832.3076

$ws.Range("H86").Value = 179.125
$ws.Range("J86").Value = 173.86957
$ws.Range("L86").Value = 521.60871
$ws.Range("N86").Value = -2893.60871

$ws.Range("H89").Value = 179.125
$ws.Range("J89").Value = 173.86957
$ws.Range("L89").Value = 1564.82613
$ws.Range("N89").Value = -13420.82613

$ws.Range("H131").Value = 1239.9403
$ws.Range("I131").Value = 2010.7693
$ws.Range("J131").Value = 1054.3704
$ws.Range("K131").Value = 6032.3079
$ws.Range("L131").Value = 3163.1112
$ws.Range("M131").Value = -992.3078999999998
$ws.Range("N131").Value = -13243.1112


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 100024
$ws.Range("J38").Value = 100024
$ws.Range("L38").Value = 100024
$ws.Range("N38").Value = -100950

$ws.Range("H102").Value = 35086.773
$ws.Range("I102").Value = 2208.65
$ws.Range("K102").Value = 2208.65
$ws.Range("M102").Value = -586.6500000000001

$ws.Range("H122").Value = 5012.225
$ws.Range("I122").Value = 4019.64
$ws.Range("J122").Value = 6666.533
$ws.Range("K122").Value = 12058.92
$ws.Range("L122").Value = 19999.599
$ws.Range("M122").Value = -9608.92
$ws.Range("N122").Value = -24899.599

$ws.Range("H132").Value = 4159.353
$ws.Range("I132").Value = 4495.8823
$ws.Range("J132").Value = 3822.8235
$ws.Range("K132").Value = 13487.6469
$ws.Range("L132").Value = 11468.4705
$ws.Range("M132").Value = -10957.6469
$ws.Range("N132").Value = -16528.4705


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2081.0908
$ws.Range("I7").Value = 1645.5294
$ws.Range("J7").Value = 3562
$ws.Range("K7").Value = 1645.5294
$ws.Range("L7").Value = 3562
$ws.Range("M7").Value = -1533.5294
$ws.Range("N7").Value = -3786

$ws.Range("H31").Value = 839.94116
$ws.Range("J31").Value = 835.7143
$ws.Range("L31").Value = 835.7143
$ws.Range("N31").Value = -1331.7143

$ws.Range("H46").Value = 1923.3077
$ws.Range("I46").Value = 474.5
$ws.Range("J46").Value = 2186.7273
$ws.Range("K46").Value = 474.5
$ws.Range("L46").Value = 2186.7273
$ws.Range("M46").Value = -286.5
$ws.Range("N46").Value = -2562.7273

$ws.Range("H126").Value = 2081.0908
$ws.Range("I126").Value = 1645.5294
$ws.Range("J126").Value = 3562
$ws.Range("K126").Value = 4936.5882
$ws.Range("L126").Value = 10686
$ws.Range("M126").Value = -2466.5882
$ws.Range("N126").Value = -15626


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 5814.8335
$ws.Range("J55").Value = 7497.25
$ws.Range("L55").Value = 7497.25
$ws.Range("N55").Value = -8051.25

$ws.Range("H113").Value = 1531.3636
$ws.Range("I113").Value = 223.66667
$ws.Range("K113").Value = 671.00001
$ws.Range("M113").Value = 1498.99999

$ws.Range("H124").Value = 22736.857
$ws.Range("J124").Value = 22736.857
$ws.Range("L124").Value = 22736.857
$ws.Range("N124").Value = -32556.857

$ws.Range("H126").Value = 22259
$ws.Range("I126").Value = 31733.242
$ws.Range("J126").Value = 2718.375
$ws.Range("K126").Value = 95199.726
$ws.Range("L126").Value = 8155.125
$ws.Range("M126").Value = -92729.726
$ws.Range("N126").Value = -13095.125

$ws.Range("H136").Value = 3113.5925
$ws.Range("I136").Value = 2803.8235
$ws.Range("J136").Value = 3640.2
$ws.Range("K136").Value = 8411.470499999999
$ws.Range("L136").Value = 10920.6
$ws.Range("M136").Value = -5861.470499999999
$ws.Range("N136").Value = -16020.6

